$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 ("RM 232") entirely - remaining rows shift up by one
$ws.Rows.Item(26).Delete()

# After the first deletion, the row that used to be 28 ("SC 92") is now row 27.
# Delete it too, so remaining rows shift up again.
$ws.Rows.Item(27).Delete()

# Now update individual F-column values (row numbers below reflect the
# post-deletion layout, matching the final state described by the diff).
$ws.Range("F3").Value = 17.64      # RM 8
$ws.Range("F5").Value = ""         # RM 14 -> now blank
$ws.Range("F21").Value = 16.58     # RM 135
$ws.Range("F23").Value = ""        # RM 140 -> now blank
$ws.Range("F32").Value = 17.39     # SC 193
